$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures in the header block ---
$ws.Range("E11").Value = 1280078      # VALOR MORA total
$ws.Range("C13").Value = 3            # Cant. Trabajadores
$ws.Range("F13").Value = 8            # Cant. Periodos

# --- Insert 6 new rows before the existing data rows to hold the new worker's periods ---
$ws.Rows("16:21").Insert()

# Copy the data-row formatting down into the freshly inserted rows (reuses the
# existing "inner row" style instead of Excel fabricating new style indices)
$ws.Range("B22:J22").Copy()
$ws.Range("B16:J21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Populate the 6 new rows for JUDITH MARGARITA DIAZ AGAMEZ ---
$periods = @("2507","2506","2505","2504","2503","2502")
for ($i = 0; $i -lt 6; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "45528478"
    $ws.Cells.Item($r, 4).Value = "JUDITH MARGARITA DIAZ AGAMEZ"
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = 203884
    $ws.Cells.Item($r, 7).Value = 5097095
}
